{"js": "// Replace each two-digit multiplication problem/answer with its updated value.\n// Old and new values are all distinct, so a straightforward sequential\n// search-and-replace over the whole document body is unambiguous and safe.\nconst replacements = [\n  [\"80\u00d744=3520\", \"72\u00d736=2592\"],\n  [\"35\u00d739=1365\", \"50\u00d776=3800\"],\n  [\"33\u00d756=1848\", \"31\u00d726=806\"],\n  [\"15\u00d724=360\", \"23\u00d715=345\"],\n  [\"85\u00d761=5185\", \"58\u00d721=1218\"],\n  [\"33\u00d795=3135\", \"16\u00d728=448\"],\n  [\"31\u00d718=558\", \"46\u00d744=2024\"],\n  [\"69\u00d742=2898\", \"57\u00d785=4845\"],\n  [\"92\u00d722=2024\", \"87\u00d722=1914\"],\n  [\"66\u00d763=4158\", \"24\u00d767=1608\"],\n  [\"84\u00d774=6216\", \"37\u00d744=1628\"],\n  [\"99\u00d759=5841\", \"55\u00d795=5225\"],\n  [\"59\u00d733=1947\", \"56\u00d783=4648\"],\n  [\"73\u00d794=6862\", \"85\u00d741=3485\"],\n  [\"79\u00d763=4977\", \"17\u00d755=935\"],\n  [\"66\u00d778=5148\", \"44\u00d753=2332\"],\n  [\"44\u00d712=528\", \"17\u00d778=1326\"],\n  [\"87\u00d740=3480\", \"67\u00d713=871\"],\n  [\"44\u00d769=3036\", \"96\u00d744=4224\"],\n  [\"96\u00d796=9216\", \"72\u00d714=1008\"],\n  [\"25\u00d725=625\", \"61\u00d738=2318\"],\n  [\"82\u00d756=4592\", \"52\u00d715=780\"],\n  [\"68\u00d761=4148\", \"28\u00d736=1008\"],\n  [\"44\u00d780=3520\", \"11\u00d733=363\"],\n  [\"42\u00d713=546\", \"81\u00d714=1134\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication problem/answer with its updated value.\n# Old and new values are all distinct, so a straightforward sequential\n# Find/Replace over the whole document body is unambiguous and safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"80\u00d744=3520\", \"72\u00d736=2592\"),\n    @(\"35\u00d739=1365\", \"50\u00d776=3800\"),\n    @(\"33\u00d756=1848\", \"31\u00d726=806\"),\n    @(\"15\u00d724=360\", \"23\u00d715=345\"),\n    @(\"85\u00d761=5185\", \"58\u00d721=1218\"),\n    @(\"33\u00d795=3135\", \"16\u00d728=448\"),\n    @(\"31\u00d718=558\", \"46\u00d744=2024\"),\n    @(\"69\u00d742=2898\", \"57\u00d785=4845\"),\n    @(\"92\u00d722=2024\", \"87\u00d722=1914\"),\n    @(\"66\u00d763=4158\", \"24\u00d767=1608\"),\n    @(\"84\u00d774=6216\", \"37\u00d744=1628\"),\n    @(\"99\u00d759=5841\", \"55\u00d795=5225\"),\n    @(\"59\u00d733=1947\", \"56\u00d783=4648\"),\n    @(\"73\u00d794=6862\", \"85\u00d741=3485\"),\n    @(\"79\u00d763=4977\", \"17\u00d755=935\"),\n    @(\"66\u00d778=5148\", \"44\u00d753=2332\"),\n    @(\"44\u00d712=528\", \"17\u00d778=1326\"),\n    @(\"87\u00d740=3480\", \"67\u00d713=871\"),\n    @(\"44\u00d769=3036\", \"96\u00d744=4224\"),\n    @(\"96\u00d796=9216\", \"72\u00d714=1008\"),\n    @(\"25\u00d725=625\", \"61\u00d738=2318\"),\n    @(\"82\u00d756=4592\", \"52\u00d715=780\"),\n    @(\"68\u00d761=4148\", \"28\u00d736=1008\"),\n    @(\"44\u00d780=3520\", \"11\u00d733=363\"),\n    @(\"42\u00d713=546\", \"81\u00d714=1134\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.Text = $oldText\n    $rng.Find.Replacement.Text = $newText\n    # wdReplaceAll = 2, wdFindContinue (Wrap) = 1\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
